$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1. Insert the new "ustadz" worksheet before "category"
# ------------------------------------------------------------------
$category = $wb.Worksheets.Item("category")
$ustadz = $wb.Worksheets.Add($category)
$ustadz.Name = "ustadz"

$ustadz.Range("A1").Value = "id"
$ustadz.Range("B1").Value = "fullname"
$ustadz.Range("C1").Value = "gender"
$ustadz.Range("D1").Value = "phone"
$ustadz.Range("E1").Value = "birthDate"

$ustadz.Range("A2").Value = 1
$ustadz.Range("B2").Value = "Ustadz Hanan Attaki"
$ustadz.Range("C2").Value = "Laki-laki"
# D2 (phone) has no data for this ustadz, but keeps the "quoted text"
# style used throughout the workbook for phone/date-like text cells.
$ustadz.Range("D2").Value = "'"
$ustadz.Range("D2").Value = ""
$ustadz.Range("E2").Value = "'1981-12-31"

# (approximate "autofit" column widths matching the saved workbook)
$ustadz.Columns.Item(1).ColumnWidth = 1.6666666666666667
$ustadz.Columns.Item(2).ColumnWidth = 17.365885416666668
$ustadz.Columns.Item(3).ColumnWidth = 5.865885416666667
$ustadz.Columns.Item(4).ColumnWidth = 12.166666666666666
$ustadz.Columns.Item(5).ColumnWidth = 9.264322916666666

# ------------------------------------------------------------------
# 2. account sheet updates
# ------------------------------------------------------------------
$account = $wb.Worksheets.Item("account")
$account.Range("E1").Value = "birthDate"
$account.Range("C2").Value = "Laki-laki"
$account.Range("L2").Value = $true
$account.Range("C3").Value = "Laki-laki"
$account.Range("G3").Value = "fulan@gmail.com"

# ------------------------------------------------------------------
# 3. content sheet updates (ustadzName column -> ustadzId column)
# ------------------------------------------------------------------
$content = $wb.Worksheets.Item("content")
$content.Range("C1").Value = "ustadzId"
$content.Range("C2").Value = 1
$content.Range("C3").Value = 1

# Column widths re-fit after the ustadzName -> ustadzId column change
$content.Columns.Item(1).ColumnWidth = 1.6666666666666667
$content.Columns.Item(3).ColumnWidth = 7.166666666666667
$content.Columns.Item(4).ColumnWidth = 15.463541666666666
$content.Columns.Item(5).ColumnWidth = 71.26432291666667
$content.Columns.Item(6).ColumnWidth = 6.065104166666667
$content.Columns.Item(7).ColumnWidth = 7.264322916666667
$content.Columns.Item(8).ColumnWidth = 6.463541666666667

# ------------------------------------------------------------------
# 4. Selections: restore each sheet's last-used cell, then leave
#    "account" as the active tab/selection (matches the workbook's
#    tabSelected flag).
# ------------------------------------------------------------------
$ustadz.Range("F5").Select() | Out-Null

$content.Activate() | Out-Null
$content.Range("E4").Select() | Out-Null

$account.Activate() | Out-Null
$account.Range("D16").Select() | Out-Null
